$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'61.779.67"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +2.30%  "
$ws.Range("D3").Value = "'3.389.55"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.34%  "
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("D5").Value = "'576.21"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.27%  "
$ws.Range("D6").Value = "'136.51"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.31%  "
$ws.Range("E7").Value = "  -0.10%  "
$ws.Range("B8").Value = "XRP"
$ws.Range("C8").Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
$ws.Range("D8").Value = "'0.475"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.13%  "
$ws.Range("B9").Value = "Toncoin"
$ws.Range("C9").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D9").Value = "'7.48"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.25%  "
$ws.Range("B10").Value = "Dogecoin"
$ws.Range("C10").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D10").Value = "'0.126"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +5.70%  "
$ws.Range("B11").Value = "Cardano"
$ws.Range("C11").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D11").Value = "'0.392"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.07%  "
$ws.Range("B12").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C12").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D12").Value = "'3.968.11"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.01%  "
$ws.Range("B13").Value = "TRON"
$ws.Range("C13").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D13").Value = "'0.122"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.07%  "
$ws.Range("B14").Value = "ShibaInu"
$ws.Range("C14").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D14").Value = "'0.0000177"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.32%  "
$ws.Range("B15").Value = "WrappedEther"
$ws.Range("C15").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D15").Value = "'3.396.48"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.21%  "
$ws.Range("B16").Value = "Avalanche"
$ws.Range("C16").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D16").Value = "'25.33"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.35%  "
$ws.Range("B17").Value = "WrappedBTC"
$ws.Range("C17").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D17").Value = "'61.796.54"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.99%  "
$ws.Range("B18").Value = "Chainlink"
$ws.Range("C18").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D18").Value = "'14.11"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +3.41%  "
$ws.Range("B19").Value = "Polkadot"
$ws.Range("C19").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D19").Value = "'5.83"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.55%  "
$ws.Range("B20").Value = "Uniswap"
$ws.Range("C20").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D20").Value = "'9.48"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.22%  "
$ws.Range("B21").Value = "BitcoinCash"
$ws.Range("C21").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D21").Value = "'391.06"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.92%  "
$ws.Range("B22").Value = "Polygon"
$ws.Range("C22").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D22").Value = "'0.566"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.99%  "
$ws.Range("B23").Value = "WrappedeETH"
$ws.Range("C23").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D23").Value = "'3.528.18"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.38%  "
$ws.Range("B24").Value = "PEPE"
$ws.Range("C24").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D24").Value = "'0.0000128"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +11.08%  "
$ws.Range("B25").Value = "Dai"
$ws.Range("C25").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D25").Value = "'1.00"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.17%  "
$ws.Range("B26").Value = "Litecoin"
$ws.Range("C26").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D26").Value = "'71.26"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.38%  "
$ws.Range("B27").Value = "Fetch.AI"
$ws.Range("C27").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D27").Value = "'1.59"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.39%  "
$ws.Range("B28").Value = "RenderToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D28").Value = "'7.63"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.25%  "
$ws.Range("B29").Value = "Binance-PegBSC-USD"
$ws.Range("C29").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D29").Value = "'0.999"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.24%  "
$ws.Range("B30").Value = "Kaspa"
$ws.Range("C30").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D30").Value = "'0.161"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +4.43%  "
$ws.Range("B31").Value = "InternetComputer(DFINITY)"
$ws.Range("C31").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D31").Value = "'8.26"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.54%  "
$ws.Range("B32").Value = "PancakeSwap"
$ws.Range("C32").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D32").Value = "'2.17"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.32%  "
$ws.Range("B33").Value = "USDe"
$ws.Range("C33").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D33").Value = "'1.00"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.11%  "
$ws.Range("B34").Value = "EthereumClassic"
$ws.Range("C34").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D34").Value = "'23.44"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.12%  "
$ws.Range("D35").Value = "'3.418.19"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.14%  "
$ws.Range("B36").Value = "NEARProtocol"
$ws.Range("C36").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D36").Value = "'5.41"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.45%  "
$ws.Range("B37").Value = "Aptos"
$ws.Range("C37").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D37").Value = "'6.92"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.24%  "
$ws.Range("B38").Value = "ImmutableX"
$ws.Range("C38").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D38").Value = "'1.55"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.95%  "
$ws.Range("B39").Value = "Monero"
$ws.Range("C39").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D39").Value = "'163.48"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.00%  "
$ws.Range("B40").Value = "Hedera"
$ws.Range("C40").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D40").Value = "'0.0788"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.63%  "
$ws.Range("B41").Value = "Stacks"
$ws.Range("C41").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D41").Value = "'1.78"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +11.98%  "
$ws.Range("B42").Value = "ONDO"
$ws.Range("C42").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D42").Value = "'1.24"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.81%  "
$ws.Range("D43").Value = "'0.783"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +4.07%  "
$ws.Range("E44").Value = "  -0.17%  "
$ws.Range("B45").Value = "Filecoin"
$ws.Range("C45").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D45").Value = "'4.44"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.40%  "
$ws.Range("B46").Value = "OKB"
$ws.Range("C46").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D46").Value = "'41.74"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.37%  "
$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").Value = "'24.38"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.81%  "
$ws.Range("B48").Value = "Cosmos"
$ws.Range("C48").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D48").Value = "'6.91"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.12%  "
$ws.Range("B49").Value = "InjectiveProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D49").Value = "'23.27"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.73%  "
$ws.Range("B50").Value = "Maker"
$ws.Range("C50").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D50").Value = "'2.345.63"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +7.77%  "
$ws.Range("B51").Value = "VeChain"
$ws.Range("C51").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D51").Value = "'0.0263"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.23%  "
